$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 22227022
$ws.Range("I62").Value = 27783276
$ws.Range("K62").Value = 27783276
$ws.Range("M62").Value = -27782652
$ws.Range("H65").Value = 22227022
$ws.Range("I65").Value = 27783276
$ws.Range("K65").Value = 138916380
$ws.Range("M65").Value = -138913260
$ws.Range("H132").Value = 7414692
$ws.Range("I132").Value = 13339722
$ws.Range("K132").Value = 40019166
$ws.Range("M132").Value = -40016636
$ws.Range("H135").Value = 40001710
$ws.Range("I135").Value = 512.7895
$ws.Range("J135").Value = 166672180
$ws.Range("K135").Value = 4615.1055
$ws.Range("L135").Value = 1500049620
$ws.Range("M135").Value = -2080.1055
$ws.Range("N135").Value = -1500054690
$ws.Range("H137").Value = 1488.6666
$ws.Range("I137").Value = 1463.8182
$ws.Range("J137").Value = 1505.75
$ws.Range("K137").Value = 4391.4546
$ws.Range("L137").Value = 4517.25
$ws.Range("M137").Value = -1841.4546
$ws.Range("N137").Value = -9617.25
$ws.Range("H141").Value = 1706.6666
$ws.Range("I141").Value = 1060
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 3180
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 2000
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1220.1666
$ws.Range("I45").Value = 1055.25
$ws.Range("K45").Value = 1055.25
$ws.Range("M45").Value = -678.25
$ws.Range("H61").Value = 90910940
$ws.Range("I61").Value = 100001736
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 100001736
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -100001524
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 2287
$ws.Range("I74").Value = 1860.4166
$ws.Range("J74").Value = 3993.3333
$ws.Range("K74").Value = 1860.4166
$ws.Range("L74").Value = 3993.3333
$ws.Range("M74").Value = -986.4166
$ws.Range("N74").Value = -5741.3333
$ws.Range("H77").Value = 2287
$ws.Range("I77").Value = 1860.4166
$ws.Range("J77").Value = 3993.3333
$ws.Range("K77").Value = 9302.083000000001
$ws.Range("L77").Value = 19966.6665
$ws.Range("M77").Value = -4934.083000000001
$ws.Range("N77").Value = -28702.6665
$ws.Range("H122").Value = 952.9167
$ws.Range("I122").Value = 960.5454999999999
$ws.Range("J122").Value = 869
$ws.Range("K122").Value = 2881.6365
$ws.Range("L122").Value = 2607
$ws.Range("M122").Value = -431.6364999999996
$ws.Range("N122").Value = -7507
$ws.Range("H132").Value = 2600.0476
$ws.Range("I132").Value = 1974.375
$ws.Range("J132").Value = 4602.2
$ws.Range("K132").Value = 5923.125
$ws.Range("L132").Value = 13806.6
$ws.Range("M132").Value = -3393.125
$ws.Range("N132").Value = -18866.6
$ws.Range("H136").Value = 90910940
$ws.Range("I136").Value = 100001736
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 300005208
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -300002658
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1289.3334
$ws.Range("I107").Value = 1080.3334
$ws.Range("K107").Value = 1080.3334
$ws.Range("M107").Value = 839.6666
$ws.Range("H110").Value = 49666.332
$ws.Range("J110").Value = 49666.332
$ws.Range("L110").Value = 49666.332
$ws.Range("N110").Value = -57846.332
$ws.Range("H130").Value = 35000
$ws.Range("J130").Value = 35000
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040
$ws.Range("H134").Value = 5662.5
$ws.Range("I134").Value = 1045.05
$ws.Range("K134").Value = 3135.15
$ws.Range("M134").Value = -600.1499999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1340.9445
$ws.Range("I31").Value = 1388
$ws.Range("J31").Value = 1293.8889
$ws.Range("K31").Value = 1388
$ws.Range("L31").Value = 1293.8889
$ws.Range("M31").Value = -1093
$ws.Range("N31").Value = -1883.8889
$ws.Range("H34").Value = 1340.9445
$ws.Range("I34").Value = 1388
$ws.Range("J34").Value = 1293.8889
$ws.Range("K34").Value = 1388
$ws.Range("L34").Value = 1293.8889
$ws.Range("M34").Value = -1186
$ws.Range("N34").Value = -1697.8889
$ws.Range("H105").Value = 1400
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -4494
$ws.Range("H109").Value = 19333.666
$ws.Range("J109").Value = 19333.666
$ws.Range("L109").Value = 19333.666
$ws.Range("N109").Value = -21413.666
$ws.Range("H132").Value = 5090.933
$ws.Range("I132").Value = 8030
$ws.Range("J132").Value = 2151.8667
$ws.Range("K132").Value = 24090
$ws.Range("L132").Value = 6455.6001
$ws.Range("M132").Value = -21560
$ws.Range("N132").Value = -11515.6001
$ws.Range("H134").Value = 33336596
$ws.Range("I134").Value = 4302.4443
$ws.Range("J134").Value = 83335030
$ws.Range("K134").Value = 12907.3329
$ws.Range("L134").Value = 250005090
$ws.Range("M134").Value = -10372.3329
$ws.Range("N134").Value = -250010160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 733480.4399999999
$ws.Range("I4").Value = 69203.766
$ws.Range("J4").Value = 2460599.8
$ws.Range("K4").Value = 207611.298
$ws.Range("L4").Value = 7381799.399999999
$ws.Range("M4").Value = -207499.298
$ws.Range("N4").Value = -7382023.399999999
$ws.Range("H131").Value = 14495434
$ws.Range("J131").Value = 3130.6206
$ws.Range("L131").Value = 9391.861800000001
$ws.Range("N131").Value = -19471.8618
$ws.Range("H132").Value = 965.5
$ws.Range("I132").Value = 959.6
$ws.Range("J132").Value = 995
$ws.Range("K132").Value = 8636.4
$ws.Range("L132").Value = 8955
$ws.Range("M132").Value = -6106.4
$ws.Range("N132").Value = -14015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 26497
$ws.Range("J86").Value = 26497
$ws.Range("L86").Value = 26497
$ws.Range("N86").Value = -28869
$ws.Range("H89").Value = 26497
$ws.Range("J89").Value = 26497
$ws.Range("L89").Value = 79491
$ws.Range("N89").Value = -91347
$ws.Range("H132").Value = 3892.7334
$ws.Range("I132").Value = 3782.75
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 11348.25
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -8818.25
$ws.Range("N132").Value = -18057.9995
$ws.Range("H135").Value = 36136.668
$ws.Range("J135").Value = 33364
$ws.Range("L135").Value = 33364
$ws.Range("N135").Value = -43504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3931
$ws.Range("I46").Value = 712.5
$ws.Range("J46").Value = 6076.6665
$ws.Range("K46").Value = 712.5
$ws.Range("L46").Value = 6076.6665
$ws.Range("M46").Value = -524.5
$ws.Range("N46").Value = -6452.6665
$ws.Range("H110").Value = 29999.5
$ws.Range("J110").Value = 29999.5
$ws.Range("L110").Value = 29999.5
$ws.Range("N110").Value = -38179.5
$ws.Range("H132").Value = 63811.168
$ws.Range("I132").Value = 19333.834
$ws.Range("J132").Value = 86049.836
$ws.Range("K132").Value = 58001.50199999999
$ws.Range("L132").Value = 258149.508
$ws.Range("M132").Value = -55471.50199999999
$ws.Range("N132").Value = -263209.508

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 9999.5
$ws.Range("J119").Value = 9999.5
$ws.Range("L119").Value = 9999.5
$ws.Range("N119").Value = -19675.5
$ws.Range("H132").Value = 2117.9
$ws.Range("I132").Value = 1936.6111
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 5809.8333
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -3279.8333
$ws.Range("N132").Value = -16308.5
$ws.Range("H136").Value = 1153.3462
$ws.Range("I136").Value = 1099.5294
$ws.Range("K136").Value = 3298.5882
$ws.Range("M136").Value = -748.5881999999997
